# Edit script: 
#  1. Insert new row 39 "Multigroup Vaccine Model" (pushes old rows 39.. down by one)
#  2. Insert new row 50 "rbranding" (pushes old rows 50.. down by one)
#  3. Rename the old "Vaccine Equity" row (now at row 58) to "Multigroup Vaccine Model"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: insert new row at 39 for "Multigroup Vaccine Model"
# ---------------------------------------------------------------------------
$ws.Rows("39").Insert()

$ws.Range("A39").Value2 = "Multigroup Vaccine Model"
$ws.Range("B39").Value2 = "A tool for exploring the effects of varying rates of vaccination among disparate socio-economic groups"
$ws.Range("C39").Value2 = "Damon Toth"
$ws.Range("D39").Value2 = "damon.toth@hcs.utah.edu"
$ws.Range("E39").Value2 = "Yes"
$ws.Range("G39").Value2 = "On development"
$ws.Range("H39").Value2 = "MIT"
$ws.Range("I39").Value2 = "R"
$ws.Range("J39").Value2 = "Modelers"
$ws.Range("K39").Value2 = "TBD"
$ws.Range("M39").Value2 = "Parameter inputs for simulating the model"
$ws.Range("N39").Value2 = "https://github.com/EpiForeSITE/multigroup-vaccine"
$ws.Range("O39").Value2 = "https://github.com/EpiForeSITE/multigroup-vaccine"

# ---------------------------------------------------------------------------
# Step 2: insert new row at 50 for "rbranding"
# ---------------------------------------------------------------------------
$ws.Rows("50").Insert()

$ws.Range("A50").Value2 = "rbranding"
$ws.Range("B50").Value2 = "An R package that facilitates the creation of Shiny applications and quarto documents with support for site-specific branding and styling."
$ws.Range("C50").Value2 = "Willy Ray"
$ws.Range("D50").Value2 = "william.ray@hsc.utah.edu"
$ws.Range("E50").Value2 = "Yes"
$ws.Range("G50").Value2 = "On development"
$ws.Range("H50").Value2 = "MIT"
$ws.Range("I50").Value2 = "R"
$ws.Range("J50").Value2 = "Public health professionals, dashboard and report builders"
$ws.Range("K50").Value2 = "TBD"
$ws.Range("L50").Value2 = "Developer Tool"
$ws.Range("N50").Value2 = "https://epiforesite.github.io/branding-package/"
$ws.Range("O50").Value2 = "https://github.com/EpiForeSITE/branding-package"

# ---------------------------------------------------------------------------
# Step 3: rename the old "Vaccine Equity" row (now shifted to row 58) to
# "Multigroup Vaccine Model". All other fields in that row stay the same.
# ---------------------------------------------------------------------------
$ws.Range("A58").Value2 = "Multigroup Vaccine Model"
